$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.959.98'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '1.676.70'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  +0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '214.84'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('E6').Value = '  -3.96%  '
$ws.Range('E7').Value = '  +0.09%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.250'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('E9').Value = '  -0.55%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '20.51'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +1.51%  '
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('D12').Value = '1.915.25'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = '1.702.17'
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('D17').Value = '26.988.97'
$ws.Range('E18').Value = '  +5.46%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '235.78'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D20').Value = '0.0₃0734'
$ws.Range('E21').Value = '  +0.06%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.44'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('E23').Value = '  -1.00%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.14'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -3.72%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '146.75'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('E26').Value = '  +0.86%  '
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('E28').Value = '  -4.25%  '
$ws.Range('E29').Value = '  +0.06%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0500'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').Value = '1.481.87'
$ws.Range('E33').Value = '  +1.73%  '
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('E35').Value = '  +4.71%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('E37').Value = '  +4.02%  '
$ws.Range('E38').Value = '  +3.11%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.908'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +1.17%  '
$ws.Range('E40').Value = '  -4.07%  '
$ws.Range('E41').Value = '  +4.51%  '
$ws.Range('E42').Value = '  +0.10%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '2.30'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.57%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '67.35'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +2.26%  '
$ws.Range('D45').Value = '1.822.58'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('E46').Value = '  +0.15%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '90.35'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('E48').Value = '  +0.87%  '
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('E50').Value = '  +2.40%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '7.77'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +1.71%  '
